$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Range("B6").Value = 254
$ws.Range("C6").Value = 49
$ws.Range("H6").Value = 152.3503937007874
$ws.Range("I6").Value = 129.265450066412
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 1.5
$ws.Range("O6").Value = 2.704724409448819
$ws.Range("P6").Value = 3.323680950706546
$ws.Range("R6").Value = 8.25
$ws.Range("S6").Value = 17.5
$ws.Range("T6").Value = 28
$ws.Range("U6").Value = 59
$ws.Range("V6").Value = 18.65748031496063
$ws.Range("W6").Value = 12.70821849597764
$ws.Range("Y6").Value = 14
$ws.Range("Z6").Value = 29
$ws.Range("AA6").Value = 39.75
$ws.Range("AC6").Value = 26.78740157480315
$ws.Range("AD6").Value = 16.95617691326643
$ws.Range("AF6").Value = 1
$ws.Range("AJ6").Value = 1.740157480314961
$ws.Range("AK6").Value = 1.083199757308116
$ws.Range("AM6").Value = 0.4210526315789473
$ws.Range("AQ6").Value = 0.5617932626842698
$ws.Range("AR6").Value = 0.3059352291189009
$ws.Range("AV6").Value = 47.55555555555556
$ws.Range("AX6").Value = 45.34535278806932
$ws.Range("AY6").Value = 16.6245696389668

# Row 7
$ws.Range("R7").Value = 7
$ws.Range("S7").Value = 38
$ws.Range("T7").Value = 56
$ws.Range("V7").Value = 48.3282208588957
$ws.Range("W7").Value = 61.82814595232369
$ws.Range("Y7").Value = 18
$ws.Range("Z7").Value = 64
$ws.Range("AA7").Value = 81.75
$ws.Range("AC7").Value = 82.88343558282209
$ws.Range("AD7").Value = 99.74264685528911
$ws.Range("AF7").Value = 1
$ws.Range("AJ7").Value = 3.285276073619632
$ws.Range("AK7").Value = 2.932297248629149

# Row 8
$ws.Range("B8").Value = 254
$ws.Range("C8").Value = 49
$ws.Range("H8").Value = 152.3503937007874
$ws.Range("I8").Value = 129.265450066412
$ws.Range("R8").Value = 30
$ws.Range("S8").Value = 46
$ws.Range("T8").Value = 60
$ws.Range("V8").Value = 61.65354330708661
$ws.Range("W8").Value = 60.31968562220263
$ws.Range("Y8").Value = 47
$ws.Range("Z8").Value = 61
$ws.Range("AA8").Value = 76
$ws.Range("AC8").Value = 90.72834645669292
$ws.Range("AD8").Value = 85.74113664665592
$ws.Range("AG8").Value = 3
$ws.Range("AH8").Value = 4
$ws.Range("AJ8").Value = 3.44488188976378
$ws.Range("AK8").Value = 2.47415761722017
$ws.Range("AS8").Value = 49
$ws.Range("AX8").Value = 152.3503937007874
$ws.Range("AY8").Value = 129.265450066412

# Row 9
$ws.Range("B9").Value = 254
$ws.Range("C9").Value = 49
$ws.Range("H9").Value = 152.3503937007874
$ws.Range("I9").Value = 129.265450066412
$ws.Range("O9").Value = 3.291338582677165
$ws.Range("P9").Value = 3.854283046930683
$ws.Range("R9").Value = 7
$ws.Range("S9").Value = 17
$ws.Range("T9").Value = 26
$ws.Range("U9").Value = 51
$ws.Range("V9").Value = 17.26771653543307
$ws.Range("W9").Value = 11.70821057058593
$ws.Range("Y9").Value = 11
$ws.Range("Z9").Value = 26
$ws.Range("AA9").Value = 37
$ws.Range("AC9").Value = 24.72834645669291
$ws.Range("AD9").Value = 15.69374084158438
$ws.Range("AF9").Value = 1
$ws.Range("AG9").Value = 2
$ws.Range("AJ9").Value = 1.598425196850394
$ws.Range("AK9").Value = 0.9679434416535295
$ws.Range("AM9").Value = 0.5
$ws.Range("AN9").Value = 0.6842105263157895
$ws.Range("AQ9").Value = 0.7274288082410751
$ws.Range("AR9").Value = 0.288957936411547
$ws.Range("AU9").Value = 31.83333333333334
$ws.Range("AV9").Value = 47.5
$ws.Range("AX9").Value = 38.01680199657272
$ws.Range("AY9").Value = 12.58556745886677
